$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

# Row 2: 'Bitcoin'
Set-TextValue $ws 'D2' '30.356.79'
Set-TextValue $ws 'E2' '  +0.59%  '

# Row 3: 'Ethereum'
Set-TextValue $ws 'D3' '1.876.40'
Set-TextValue $ws 'E3' '  +0.50%  '

# Row 5: 'BNB'
Set-TextValue $ws 'D5' '243.81'
Set-TextValue $ws 'E5' '  +3.95%  '

# Row 8: 'Cardano'
Set-TextValue $ws 'D8' '0.2879'
Set-TextValue $ws 'E8' '  +0.93%  '

# Row 9: 'Dogecoin'
Set-TextValue $ws 'D9' '0.06518'
Set-TextValue $ws 'E9' '  -0.71%  '

# Row 10: 'Solana'
Set-TextValue $ws 'D10' '21.28'
Set-TextValue $ws 'E10' '  -0.59%  '

# Row 11: 'TRON'
Set-TextValue $ws 'E11' '  -0.06%  '

# Row 12: 'WrappedEther'
Set-TextValue $ws 'D12' '1.876.89'
Set-TextValue $ws 'E12' '  +0.24%  '

# Row 13: 'Polygon' -> 'Litecoin'
Set-TextValue $ws 'B13' 'Litecoin'
Set-TextValue $ws 'C13' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws 'D13' '96.19'
Set-TextValue $ws 'E13' '  +0.20%  '

# Row 14: 'Litecoin' -> 'Polygon'
Set-TextValue $ws 'B14' 'Polygon'
Set-TextValue $ws 'C14' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws 'D14' '0.7332'
Set-TextValue $ws 'E14' '  +5.97%  '

# Row 15: 'Polkadot'
Set-TextValue $ws 'D15' '5.121'
Set-TextValue $ws 'E15' '  +0.30%  '

# Row 16: 'BitcoinCash'
Set-TextValue $ws 'D16' '275.17'
Set-TextValue $ws 'E16' '  +3.31%  '

# Row 17: 'WrappedBTC'
Set-TextValue $ws 'D17' '30.349.08'
Set-TextValue $ws 'E17' '  +0.59%  '

# Row 18: 'Avalanche'
Set-TextValue $ws 'D18' '13.36'
Set-TextValue $ws 'E18' '  -2.25%  '

# Row 19: 'ShibaInu'
Set-TextValue $ws 'D19' '0.000007527'
Set-TextValue $ws 'E19' '  -2.51%  '

# Row 20: 'Dai'
Set-TextValue $ws 'E20' '  +0.13%  '

# Row 21: 'WrappedliquidstakedEther2.0'
Set-TextValue $ws 'D21' '2.125.62'
Set-TextValue $ws 'E21' '  +0.22%  '

# Row 22: 'BinanceUSD'
Set-TextValue $ws 'E22' '  +0.17%  '

# Row 23: 'Uniswap'
Set-TextValue $ws 'D23' '5.220'
Set-TextValue $ws 'E23' '  -0.69%  '

# Row 24: 'Chainlink'
Set-TextValue $ws 'D24' '6.154'
Set-TextValue $ws 'E24' '  -0.19%  '

# Row 25: 'Cosmos'
Set-TextValue $ws 'D25' '9.211'
Set-TextValue $ws 'E25' '  -2.77%  '

# Row 26: 'Monero'
Set-TextValue $ws 'D26' '164.20'
Set-TextValue $ws 'E26' '  -1.41%  '

# Row 27: 'EthereumClassic'
Set-TextValue $ws 'D27' '18.91'
Set-TextValue $ws 'E27' '  +1.20%  '

# Row 28: 'LidoDAOToken'
Set-TextValue $ws 'D28' '1.949'
Set-TextValue $ws 'E28' '  +0.66%  '

# Row 30: 'Stellar'
Set-TextValue $ws 'D30' '0.09955'
Set-TextValue $ws 'E30' '  +0.19%  '

# Row 31: 'PancakeSwap'
Set-TextValue $ws 'D31' '1.507'
Set-TextValue $ws 'E31' '  +3.26%  '

# Row 32: 'Filecoin'
Set-TextValue $ws 'D32' '4.302'
Set-TextValue $ws 'E32' '  -1.36%  '

# Row 33: 'InternetComputer(DFINITY)'
Set-TextValue $ws 'D33' '4.075'
Set-TextValue $ws 'E33' '  +0.74%  '

# Row 34: 'Hedera'
Set-TextValue $ws 'D34' '0.04740'
Set-TextValue $ws 'E34' '  +0.35%  '

# Row 35: 'ARBITRUM'
Set-TextValue $ws 'D35' '1.120'
Set-TextValue $ws 'E35' '  -0.91%  '

# Row 36: 'ImmutableX'
Set-TextValue $ws 'D36' '0.6928'
Set-TextValue $ws 'E36' '  -1.06%  '

# Row 37: 'HuobiToken'
Set-TextValue $ws 'E37' '  +0.05%  '

# Row 38: 'VeChain'
Set-TextValue $ws 'D38' '0.01852'
Set-TextValue $ws 'E38' '  -0.24%  '

# Row 39: 'MXToken'
Set-TextValue $ws 'D39' '2.743'
Set-TextValue $ws 'E39' '  -0.87%  '

# Row 40: 'FraxShare'
Set-TextValue $ws 'D40' '6.282'
Set-TextValue $ws 'E40' '  +0.10%  '

# Row 41: 'TrustWalletToken'
Set-TextValue $ws 'D41' '0.8421'
Set-TextValue $ws 'E41' '  +0.65%  '

# Row 42: 'TheSandbox' -> 'PaxDollar'
Set-TextValue $ws 'B42' 'PaxDollar'
Set-TextValue $ws 'C42' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws 'D42' '1.001'
Set-TextValue $ws 'E42' '  +0.11%  '

# Row 43: 'RenderToken' -> 'TheSandbox'
Set-TextValue $ws 'B43' 'TheSandbox'
Set-TextValue $ws 'C43' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws 'D43' '0.4161'
Set-TextValue $ws 'E43' '  +0.49%  '

# Row 44: 'PaxDollar' -> 'RenderToken'
Set-TextValue $ws 'B44' 'RenderToken'
Set-TextValue $ws 'C44' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws 'D44' '1.904'
Set-TextValue $ws 'E44' '  -1.53%  '

# Row 45: 'Aave'
Set-TextValue $ws 'D45' '69.07'
Set-TextValue $ws 'E45' '  -4.68%  '

# Row 46: 'Quant'
Set-TextValue $ws 'D46' '101.77'
Set-TextValue $ws 'E46' '  -1.07%  '

# Row 47: 'EnergySwap'
Set-TextValue $ws 'D47' '9.268'
Set-TextValue $ws 'E47' '  +1.56%  '

# Row 48: 'Aptos'
Set-TextValue $ws 'D48' '7.072'
Set-TextValue $ws 'E48' '  -0.37%  '

# Row 49: 'Elrond'
Set-TextValue $ws 'D49' '35.11'
Set-TextValue $ws 'E49' '  +1.56%  '

# Row 50: 'Maker'
Set-TextValue $ws 'D50' '911.11'
Set-TextValue $ws 'E50' '  -6.07%  '

# Row 51: 'Cronos'
Set-TextValue $ws 'D51' '0.05587'
Set-TextValue $ws 'E51' '  -1.03%  '
